# pbmc_01_all.xlsx update:
# The marker-gene count table lost three gene columns (CD2, CD14, NKG7),
# keeping CD3D, S100A8, GZMA and CD79A (in that order), and the CD3D count
# for Cell_3 was corrected from 13 to 10. The sheet view was also reset to
# a plain 100% zoom with Column A narrowed to fit the Cell_N labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the NKG7 (G), CD14 (D) and CD2 (B) columns entirely. Deleting from
# the right-most column inward keeps the remaining column letters stable
# between calls, leaving A, CD3D, S100A8, GZMA, CD79A.
$ws.Range("G1:G5").EntireColumn.Delete()
$ws.Range("D1:D5").EntireColumn.Delete()
$ws.Range("B1:B5").EntireColumn.Delete()

# Correct Cell_3's CD3D count (now column B, row 4).
$ws.Range("B4").Value = 10

# Reset the view: normal 100% zoom and no stale selection past the data.
$excel.ActiveWindow.Zoom = 100
$ws.Range("A1").Select()

# Narrow column A to fit the short Cell_1..Cell_4 labels.
$ws.Columns.Item(1).ColumnWidth = 5
